$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the existing row 4 (old row 4 content moves down to
# row 5, and a brand-new row 6 is appended) before rewriting row 4 in place.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# The tcid column (A) stores its numeric-looking ids as text, matching the
# rest of the sheet (and the numberStoredAsText ignoredErrors rule below).
$ws.Range("A4:A6").NumberFormat = "@"

# Row 4 is rewritten with the new testcase (140741 - prepopulated data check).
$ws.Range("A4").Value = "140741"
$ws.Range("B4").Value = 'Client->Create Client->To verify that the user can view the "Prepopulated" data in the client information screen.'
$ws.Range("C4").Value = "type,type,type,dropdown,type,type,type,type,click"
$ws.Range("D4").Value = "faker,faker,faker,CA,faker,faker,faker,faker,no value"
$ws.Range("E4").Value = "client_name,street,city,state,zip,mobile_number,ext,website,save"
$ws.Range("F4").Value = "yes"

# Row 5 is the original 140688 testcase, with its description text tweaked
# (Client information -> General information).
$ws.Range("A5").Value = "140688"
$ws.Range("B5").Value = 'Client->Create Client->General Information->To verify if the user is able to click on the "Next" button and navigate to the General information screen.'
$ws.Range("C5").Value = "click,assert"
$ws.Range("D5").Value = "no value, 2.General Information "
$ws.Range("E5").Value = "next,next_tab"
$ws.Range("F5").Value = "yes"

# Row 6 is a brand-new testcase (140728 - finish button / summary screen).
$ws.Range("A6").Value = "140728"
$ws.Range("B6").Value = "Client->Create Client->To verify that the user is able to navigate to the Client Summary screen after clicking the FINISH button."
$ws.Range("C6").Value = "click,assert"
$ws.Range("D6").Value = "no value,Summary"
$ws.Range("E6").Value = "finish,summary_contains"
$ws.Range("F6").Value = "yes"
